$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "257.58"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "4.86%"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-3.89%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.221"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.32%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05922"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "3.72%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.664"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.59%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8630"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.00%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.023"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "15.07%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1415"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "2.02%"

$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07182"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.13%"

$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.03135"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.70%"

$ws.Range("B12").Value = "BitMartToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09237"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.18%"

$ws.Range("B13").Value = "BitForexToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.001545"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "1.11%"

$ws.Range("B14").Value = "One"
$ws.Range("C14").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0006077"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "1.49%"

$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005681"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-4.62%"

$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.497"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.04%"

$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.268"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.76%"

$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.192"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.95%"

$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3143"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.73%"

$ws.Range("B20").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C20").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.03558"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "6.69%"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.23%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.519"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.60%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04178"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "2.62%"

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.45%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001219"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.28%"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "8.67%"

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-0.01%"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001484"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "2.73%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03818"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.005647"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "51.15%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1103"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "3.37%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002199"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-10.58%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01082"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "14.71%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005428"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "2.92%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000750"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.00%"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "22.45%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002235"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-0.96%"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002099"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.00%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0001999"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.00%"

